# Generate Report for Handoff
# Rows 4-7 in the "zh-cn" and "de-de" localization-status sheets represented
# files that hadn't been handed off yet (Priority "low", no handback info
# written). A handoff run completed for those files, so:
#   - Priority flips from "low" to "ht"
#   - The "Latest Handoff Datetime" column is refreshed with the new
#     handoff timestamp produced by that run.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-28 18:31:29"
}

$dede = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-28 18:31:35"
}
